# Update countries & provincias Spain
# Re-applies the 06:16 COVID-19 snapshot refresh:
#  - Australia, Kazajistan and Polinesia Francesa received updated case
#    counts, which moved them up past their (alphabetically-unrelated,
#    count-sorted) neighbours in the country list; everything below
#    cascades down by one slot until it reaches a row whose country
#    already sits where it should.
#  - India's counts were refreshed in place (no reordering).
#  - The "last updated" footer timestamp moved from 05:46 to 06:16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $total, $new, $active, $recovered, $critical, $deathsToday, $deaths) {
    $ws.Range("A$row").Value = $country
    $ws.Range("B$row").Value = $total
    $ws.Range("C$row").Value = $new
    $ws.Range("D$row").Value = $active
    $ws.Range("E$row").Value = $recovered
    $ws.Range("F$row").Value = $critical
    $ws.Range("G$row").Value = $deathsToday
    $ws.Range("H$row").Value = $deaths
}

# Australia jumps ahead of Brasil with refreshed numbers; Brasil drops
# one row, keeping its previous (unchanged) figures.
Set-Row 21 "Australia" 1973 86 118 1847 11 1 8
Set-Row 22 "Brasil"    1924  0   2 1888 18 0 34

# India: same position, refreshed figures only.
$ws.Range("B44").Value = 508
$ws.Range("C44").Value = 9
$ws.Range("E44").Value = 461

# Kazajistan jumps ahead of Oman with refreshed numbers; Oman drops one
# row, keeping its previous (unchanged) figures.
Set-Row 100 "Kazajistan" 66 4  0 66 0 0 0
Set-Row 101 "Oman"       66 0 17 49 0 0 0

# Polinesia Francesa jumps ahead of Monaco/Guyana/Guayana Francesa/
# Jamaica/Togo with refreshed numbers; each of those five countries
# cascades down one row, keeping its own previous figures.
Set-Row 127 "Polinesia Francesa" 23 5  0 23 0 0 0
Set-Row 128 "Monaco"             23 0  1 22 0 0 0
Set-Row 130 "Guyana"             20 0  0 19 0 0 1
Set-Row 131 "Guayana Francesa"   20 0  6 14 0 0 0
Set-Row 132 "Jamaica"            19 0  2 16 0 0 0
Set-Row 133 "Togo"               18 0  0 18 0 0 0

# Footer timestamp.
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 06:16"
